$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 ("Create UI Framework") currently has an empty B column (AC 1).
# Add the missing acceptance criteria text.
$ws.Range("B8").Value = "UI must be capable of testing implemented functions"
